$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new CMIP6 Specialization ID values for key_properties rows.
# Order matters here: it reproduces the exact shared-string table ordering
# from the original author's edit (values were not entered strictly top
# to bottom - e.g. C6 before C5, and C32 before C31).
$ws.Range("C3").Value  = "cmip6.ocean.key_properties.bathymetry.type"
$ws.Range("C4").Value  = "cmip6.ocean.key_properties.bathymetry.reference_dates"
$ws.Range("C6").Value  = "cmip6.ocean.grid.discretisation.horizontal.scheme"
$ws.Range("C5").Value  = "cmip6.ocean.grid.discretisation.horizontal.type"
$ws.Range("C8").Value  = "cmip6.ocean.key_properties.prognostic_variables"
$ws.Range("C9").Value  = "cmip6.ocean.key_properties.model_family"
$ws.Range("C10").Value = "cmip6.ocean.key_properties.nonoceanic_waters.isolated_seas"
$ws.Range("C11").Value = "cmip6.ocean.key_properties.nonoceanic_waters.river_mouth"
$ws.Range("C25").Value = "cmip6.ocean.key_properties.seawater_properties.eos_type"
$ws.Range("C26").Value = "cmip6.ocean.key_properties.seawater_properties.ocean_freezing_point"
$ws.Range("C28").Value = "cmip6.ocean.key_properties.seawater_properties.ocean_specific_heat"
$ws.Range("C30").Value = "cmip6.ocean.timestepping_framework.barotropic_momentum_scheme.barotropic_momentum"
$ws.Range("C32").Value = "cmip6.ocean.timestepping_framework.barotropic_solver_scheme.type"
$ws.Range("C31").Value = "cmip6.ocean.timestepping_framework.barotropic_solver_scheme.barotropic_solver"
$ws.Range("C33").Value = "cmip6.ocean.timestepping_framework.timestepping_attributes.time_step"
$ws.Range("C34").Value = "cmip6.ocean.timestepping_framework.timestepping_tracers_scheme.tracers"

# Restore the selection and adjust the view (scrolled + zoomed to 150%)
$ws.Range("C34").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 150
